$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -3.40%  "
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").Value = "  -7.46%  "
$ws.Range("E8").Value = "  +14.72%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("E11").Value = "  +8.20%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  -7.45%  "
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  -9.38%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  -4.14%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +6.58%  "
$ws.Range("E31").Value = "  +7.26%  "
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("E34").Value = "  -12.76%  "
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("E36").Value = "  -4.06%  "
$ws.Range("E37").Value = "  -5.61%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  -6.87%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("E47").Value = "  +13.04%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("E49").Value = "  -7.25%  "
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  -6.27%  "
